$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price observation was recorded for "Ajo" (Chino / Primera) and
# inserted into the dataset right before the existing 2021-09-28 row, pushing
# every following row down by one (the old last row, 972, becomes 973).
$ws.Rows(866).Insert()

$ws.Cells.Item(866, 1).Value = 6
$ws.Cells.Item(866, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(866, 3).Value = "Metropolitana"
$ws.Cells.Item(866, 4).Value = 44946
$ws.Cells.Item(866, 5).Value = 13
$ws.Cells.Item(866, 6).Value = 100112003
$ws.Cells.Item(866, 7).Value = "Ajo"
$ws.Cells.Item(866, 8).Value = "Chino"
$ws.Cells.Item(866, 9).Value = "Primera"
$ws.Cells.Item(866, 10).Value = 1500
$ws.Cells.Item(866, 11).Value = 14000
$ws.Cells.Item(866, 12).Value = 15000
$ws.Cells.Item(866, 13).Value = 14400
$ws.Cells.Item(866, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(866, 15).Value = "China"
$ws.Cells.Item(866, 16).Value = 1440
$ws.Cells.Item(866, 17).Value = 10
$ws.Cells.Item(866, 18).Value = "Hortaliza"
